# Daily attendance processing - 2025-11-09 05:21:22
#
# Normalises the "Recorded By" column (G) so that automated/system
# recorders are listed first in the comma-separated attendee list:
#   - if "System" is present, move it to the front of the list
#   - else if "admin@admin.com" is present, move it to the front
# Single-value cells (no comma) and cells that already start with the
# priority entry are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reorder-RecordedBy {
    param([string]$text)

    $parts = @($text -split ", ")
    if ($parts.Count -le 1) {
        return $text
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }

    $hasAdmin = $false
    foreach ($p in $parts) {
        if ($p.Equals("admin@admin.com")) { $hasAdmin = $true }
    }

    if ($hasSystem) {
        $rest = @($parts | Where-Object { -not $_.Equals("System") })
        $newParts = @("System") + $rest
        return ($newParts -join ", ")
    }
    elseif ($hasAdmin) {
        $rest = @($parts | Where-Object { -not $_.Equals("admin@admin.com") })
        $newParts = @("admin@admin.com") + $rest
        return ($newParts -join ", ")
    }

    return $text
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($current -ne $null) {
        $updated = Reorder-RecordedBy $current
        if (-not $updated.Equals($current)) {
            $cell.Value = $updated
        }
    }
}
